$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Admin")

# New training session row (row 24), formatted like the row above it
# (row 23): A:E carry data, no entries in F/G (no links yet).
$ws.Range("A23:E23").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A24").Value = "Oct 28, 2024"
$ws.Range("B24").Value = "09:00 AM - 11:00 PM"
$ws.Range("C24").Value = "3 hours"
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = "Completed"

$ws.Range("E24").Select()

$wb.Save()
